# Auto-generated Excel COM-interop script to apply numeric updates
# described by the authoritative diff (Brynhildr_Profits workbook).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 477.2
$ws.Range("I28").Value = 520.75
$ws.Range("K28").Value = 520.75
$ws.Range("M28").Value = -35.75
$ws.Range("H41").Value = 473.52173
$ws.Range("I41").Value = 295.69232
$ws.Range("J41").Value = 704.7
$ws.Range("K41").Value = 295.69232
$ws.Range("L41").Value = 704.7
$ws.Range("M41").Value = 144.30768
$ws.Range("N41").Value = -1584.7
$ws.Range("H98").Value = 2066.8235
$ws.Range("I98").Value = 2188.5334
$ws.Range("K98").Value = 2188.5334
$ws.Range("M98").Value = -690.5333999999998
$ws.Range("H107").Value = 7994.1055
$ws.Range("I107").Value = 6868.625
$ws.Range("K107").Value = 6868.625
$ws.Range("M107").Value = -4948.625
$ws.Range("H122").Value = 2066.8235
$ws.Range("I122").Value = 2188.5334
$ws.Range("K122").Value = 6565.600199999999
$ws.Range("M122").Value = -4115.600199999999
$ws.Range("H129").Value = 317922.28
$ws.Range("I129").Value = 367755.38
$ws.Range("K129").Value = 1103266.14
$ws.Range("M129").Value = -1098266.14
$ws.Range("H132").Value = 1692.7347
$ws.Range("I132").Value = 1568.4651
$ws.Range("J132").Value = 2583.3333
$ws.Range("K132").Value = 4705.3953
$ws.Range("L132").Value = 7749.999899999999
$ws.Range("M132").Value = -2175.3953
$ws.Range("N132").Value = -12809.9999
$ws.Range("H135").Value = 4491.864
$ws.Range("I135").Value = 398.44446
$ws.Range("K135").Value = 3586.00014
$ws.Range("M135").Value = -1051.00014

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 933.73334
$ws.Range("J2").Value = 1500
$ws.Range("L2").Value = 1500
$ws.Range("N2").Value = -1726
$ws.Range("H32").Value = 9428.171
$ws.Range("I32").Value = 4171.6597
$ws.Range("K32").Value = 4171.6597
$ws.Range("M32").Value = -3884.6597
$ws.Range("H61").Value = 1432006.6
$ws.Range("I61").Value = 3479.3962
$ws.Range("K61").Value = 3479.3962
$ws.Range("M61").Value = -3267.3962
$ws.Range("H74").Value = 1691448.4
$ws.Range("I74").Value = 2226504.5
$ws.Range("K74").Value = 2226504.5
$ws.Range("M74").Value = -2225630.5
$ws.Range("H77").Value = 1691448.4
$ws.Range("I77").Value = 2226504.5
$ws.Range("K77").Value = 11132522.5
$ws.Range("M77").Value = -11128154.5
$ws.Range("H102").Value = 2899.9
$ws.Range("I102").Value = 2999.889
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 2999.889
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -1377.889
$ws.Range("N102").Value = -5244
$ws.Range("H116").Value = 933.73334
$ws.Range("J116").Value = 1500
$ws.Range("L116").Value = 1500
$ws.Range("N116").Value = -6088
$ws.Range("H132").Value = 965807.9399999999
$ws.Range("I132").Value = 1140817.9
$ws.Range("J132").Value = 3253.5
$ws.Range("K132").Value = 3422453.7
$ws.Range("L132").Value = 9760.5
$ws.Range("M132").Value = -3419923.7
$ws.Range("N132").Value = -14820.5
$ws.Range("H134").Value = 58606
$ws.Range("J134").Value = 58606
$ws.Range("L134").Value = 58606
$ws.Range("N134").Value = -68746
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("H136").Value = 1432006.6
$ws.Range("I136").Value = 3479.3962
$ws.Range("K136").Value = 10438.1886
$ws.Range("M136").Value = -7888.188600000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 79750
$ws.Range("J2").Value = 79750
$ws.Range("L2").Value = 79750
$ws.Range("N2").Value = -79976
$ws.Range("H3").Value = 933.73334
$ws.Range("J3").Value = 1500
$ws.Range("L3").Value = 1500
$ws.Range("N3").Value = -1728
$ws.Range("H20").Value = 45458
$ws.Range("I20").Value = 63420.59
$ws.Range("K20").Value = 63420.59
$ws.Range("M20").Value = -63173.59
$ws.Range("H105").Value = 6116.143
$ws.Range("I105").Value = 5060.6665
$ws.Range("J105").Value = 8016
$ws.Range("K105").Value = 5060.6665
$ws.Range("L105").Value = 8016
$ws.Range("M105").Value = -3313.6665
$ws.Range("N105").Value = -11510
$ws.Range("H134").Value = 5755963
$ws.Range("I134").Value = 5625
$ws.Range("K134").Value = 16875
$ws.Range("M134").Value = -14340

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 6586
$ws.Range("I68").Value = 497.16666
$ws.Range("J68").Value = 7727.6562
$ws.Range("K68").Value = 1491.49998
$ws.Range("L68").Value = 23182.9686
$ws.Range("M68").Value = -680.4999800000001
$ws.Range("N68").Value = -24804.9686
$ws.Range("H71").Value = 6586
$ws.Range("I71").Value = 497.16666
$ws.Range("J71").Value = 7727.6562
$ws.Range("K71").Value = 4474.49994
$ws.Range("L71").Value = 69548.90580000001
$ws.Range("M71").Value = -418.4999399999997
$ws.Range("N71").Value = -77660.90580000001
$ws.Range("H131").Value = 5295032
$ws.Range("I131").Value = 1333.8182
$ws.Range("J131").Value = 11118100
$ws.Range("K131").Value = 4001.4546
$ws.Range("L131").Value = 33354300
$ws.Range("M131").Value = 1038.5454
$ws.Range("N131").Value = -33364380

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3678
$ws.Range("I102").Value = 3678
$ws.Range("K102").Value = 3678
$ws.Range("M102").Value = -2056
$ws.Range("H122").Value = 6179.6294
$ws.Range("I122").Value = 6234
$ws.Range("J122").Value = 5500
$ws.Range("K122").Value = 18702
$ws.Range("L122").Value = 16500
$ws.Range("M122").Value = -16252
$ws.Range("N122").Value = -21400
$ws.Range("H126").Value = 6098.5557
$ws.Range("J126").Value = 3133.3333
$ws.Range("L126").Value = 9399.999899999999
$ws.Range("N126").Value = -14339.9999
$ws.Range("H132").Value = 10502.652
$ws.Range("I132").Value = 6550.524
$ws.Range("J132").Value = 52000
$ws.Range("K132").Value = 19651.572
$ws.Range("L132").Value = 156000
$ws.Range("M132").Value = -17121.572
$ws.Range("N132").Value = -161060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1861.3914
$ws.Range("I22").Value = 682.25
$ws.Range("J22").Value = 2109.6316
$ws.Range("K22").Value = 682.25
$ws.Range("L22").Value = 2109.6316
$ws.Range("M22").Value = -387.25
$ws.Range("N22").Value = -2699.6316
$ws.Range("H27").Value = 1861.3914
$ws.Range("I27").Value = 682.25
$ws.Range("J27").Value = 2109.6316
$ws.Range("K27").Value = 682.25
$ws.Range("L27").Value = 2118.0527
$ws.Range("M27").Value = -575.25
$ws.Range("N27").Value = -2323.6316
$ws.Range("H40").Value = 6283.4287
$ws.Range("I40").Value = 6747.5
$ws.Range("J40").Value = 5664.6665
$ws.Range("K40").Value = 6747.5
$ws.Range("L40").Value = 5664.6665
$ws.Range("M40").Value = -6611.5
$ws.Range("N40").Value = -5936.6665
$ws.Range("H55").Value = 1623.9546
$ws.Range("J55").Value = 1425.0769
$ws.Range("L55").Value = 1425.0769
$ws.Range("N55").Value = -1771.0769
$ws.Range("H61").Value = 10123.571
$ws.Range("I61").Value = 9463.723
$ws.Range("K61").Value = 9463.723
$ws.Range("M61").Value = -9261.723
$ws.Range("H111").Value = 52500
$ws.Range("J111").Value = 52500
$ws.Range("L111").Value = 52500
$ws.Range("H113").Value = 10123.571
$ws.Range("I113").Value = 9463.723
$ws.Range("K113").Value = 9463.723
$ws.Range("M113").Value = -7293.723
$ws.Range("H122").Value = 3284.2856
$ws.Range("J122").Value = 4245
$ws.Range("L122").Value = 12735
$ws.Range("N122").Value = -17635
$ws.Range("H134").Value = 59000
$ws.Range("J134").Value = 59000
$ws.Range("L134").Value = 59000
$ws.Range("N134").Value = -69140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7706.7856
$ws.Range("J62").Value = 8884.556
$ws.Range("L62").Value = 8884.556
$ws.Range("N62").Value = -10132.556
$ws.Range("H65").Value = 7706.7856
$ws.Range("J65").Value = 8884.556
$ws.Range("L65").Value = 44422.78
$ws.Range("N65").Value = -50662.78
$ws.Range("H122").Value = 44714.42
$ws.Range("I122").Value = 1279.8422
$ws.Range("J122").Value = 162608.28
$ws.Range("K122").Value = 3839.5266
$ws.Range("L122").Value = 487824.84
$ws.Range("M122").Value = -1389.5266
$ws.Range("N122").Value = -492724.84
$ws.Range("I136").Value = 3346065.8
$ws.Range("J136").Value = 50002750
$ws.Range("K136").Value = 10038197.4
$ws.Range("L136").Value = 150008250
$ws.Range("M136").Value = -10035647.4
$ws.Range("N136").Value = -150013350

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N111").Value = -60680

Write-Output "Applied numeric updates to Brynhildr_Profits sheets."